# Insert a new data row before row 263 (shifts rows 263:377 down to 264:378)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(263).Insert()

$ws.Range("A263").Value = 10
$ws.Range("B263").Value = 'Vega Modelo de Temuco'
$ws.Range("C263").Value = 'La Araucanía'
$ws.Range("D263").Value = 44755
$ws.Range("E263").Value = 9
$ws.Range("F263").Value = 100112040
$ws.Range("G263").Value = 'Cilantro'
$ws.Range("H263").Value = 'Sin especificar'
$ws.Range("I263").Value = 'Primera'
$ws.Range("J263").Value = 80
$ws.Range("K263").Value = 4600
$ws.Range("L263").Value = 4600
$ws.Range("M263").Value = 4600
$ws.Range("N263").Value = '$/docena de atados (2 kilos)'
$ws.Range("O263").Value = 'Región Metropolitana'
$ws.Range("P263").Value = 2300
$ws.Range("Q263").Value = 2
$ws.Range("R263").Value = 'Hortaliza'
